$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (preserve rich-text run formatting via Characters) ---

# A8: "Volume 31   Number  37" -> "...38" (replace just the trailing issue number run)
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "38"

# C9: "Report Covering the Week  9/9/2024  Through  9/15/2024"
#     -> "...9/16/2024  Through  9/22/2024" (replace later substring first so earlier offsets stay valid)
$c9 = $ws.Range("C9")
$c9.Characters(46, 9).Text = "9/22/2024"
$c9.Characters(27, 8).Text = "9/16/2024"

# --- Crime-stat table updates (rows 15-30) ---

# Row 15
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -50

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 60
$ws.Range("J16").Value = 87
$ws.Range("K16").Value = -31.034482758620
$ws.Range("L16").Value = -26.829268292682
$ws.Range("M16").Value = -63.414634146341
$ws.Range("N16").Value = -87.179487179487

# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 40
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -40.909090909090
$ws.Range("I17").Value = 154
$ws.Range("J17").Value = 176
$ws.Range("K17").Value = -12.5
$ws.Range("L17").Value = -3.75
$ws.Range("M17").Value = 22.222222222222
$ws.Range("N17").Value = -20.207253886010

# Row 18
$ws.Range("D18").Value = 1
$ws.Range("D18").NumberFormat = '#,##0'
$ws.Range("E18").Value = 0
$ws.Range("E18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = -50
$ws.Range("J18").Value = 39
$ws.Range("K18").Value = -10.256410256410
$ws.Range("L18").Value = -51.388888888888
$ws.Range("M18").Value = -80.978260869565
$ws.Range("N18").Value = -92.440604751619

# Row 19
$ws.Range("F19").Value = 22
$ws.Range("G19").Value = 14
$ws.Range("H19").Value = 57.142857142857
$ws.Range("I19").Value = 211
$ws.Range("J19").Value = 156
$ws.Range("K19").Value = 35.256410256410
$ws.Range("L19").Value = 17.877094972067
$ws.Range("M19").Value = 26.347305389221
$ws.Range("N19").Value = -27.491408934707

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 44.444444444444
$ws.Range("I20").Value = 106
$ws.Range("J20").Value = 111
$ws.Range("K20").Value = -4.504504504504
$ws.Range("L20").Value = -5.357142857142
$ws.Range("M20").Value = 29.268292682926
$ws.Range("N20").Value = -93.116883116883

# Row 21
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 11
$ws.Range("E21").Value = 54.545454545454
$ws.Range("F21").Value = 61
$ws.Range("G21").Value = 56
$ws.Range("H21").Value = 8.928571428571
$ws.Range("I21").Value = 579
$ws.Range("J21").Value = 585
$ws.Range("K21").Value = -1.025641025641
$ws.Range("L21").Value = -6.763285024154
$ws.Range("M21").Value = -22.072678331090
$ws.Range("N21").Value = -80.570469798657

# Row 22
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G22").Value = 1
$ws.Range("G22").NumberFormat = '#,##0'
$ws.Range("H22").Value = -100
$ws.Range("H22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J22").Value = 8
$ws.Range("K22").Value = -12.5

# Row 23
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 2
$ws.Range("D23").NumberFormat = '#,##0'
$ws.Range("E23").Value = 50
$ws.Range("E23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 74
$ws.Range("J23").Value = 90
$ws.Range("K23").Value = -17.777777777777
$ws.Range("L23").Value = -2.631578947368
$ws.Range("M23").Value = 76.190476190476

# Row 24
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = -25
$ws.Range("F24").Value = 33
$ws.Range("G24").Value = 32
$ws.Range("H24").Value = 3.125
$ws.Range("I24").Value = 382
$ws.Range("J24").Value = 420
$ws.Range("K24").Value = -9.047619047619
$ws.Range("L24").Value = -5.445544554455
$ws.Range("M24").Value = -5.679012345679

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 50
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = -14.285714285714
$ws.Range("I25").Value = 82
$ws.Range("J25").Value = 110
$ws.Range("K25").Value = -25.454545454545
$ws.Range("L25").Value = -46.052631578947

# Row 26
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -40
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = -3.703703703703
$ws.Range("I26").Value = 255
$ws.Range("J26").Value = 271
$ws.Range("K26").Value = -5.904059040590
$ws.Range("L26").Value = -12.371134020618
$ws.Range("M26").Value = -27.556818181818

# Row 27
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = -50

# Row 28
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -25
$ws.Range("J28").Value = 23
$ws.Range("K28").Value = 0

# Row 29
$ws.Range("M29").Value = -43.478260869565

# Row 30
$ws.Range("M30").Value = -40

